$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 887-888 (pushes existing 887.. down to 889..)
$ws.Rows("887:888").Insert()

# Row 887: 2026/02/27, 金, 19, 201
$ws.Range("A887").Value = "'2026/02/27"
$ws.Range("A887").ClearFormats()
$ws.Range("B887").Value = "金"
$ws.Range("C887").Value = 19
$ws.Range("D887").Value = 201

# Row 888: 2026/02/27, 金, 22, 201
$ws.Range("A888").Value = "'2026/02/27"
$ws.Range("A888").ClearFormats()
$ws.Range("B888").Value = "金"
$ws.Range("C888").Value = 22
$ws.Range("D888").Value = 201

Write-Host "Inserted rows 887-888"
